# Update "想去人数" (want-to-go count) values in both the "展览" sheet
# and the "全部类型" sheet (which duplicates the same rows).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 579
    $ws.Range("F3").Value = 55
    $ws.Range("F4").Value = 28
    $ws.Range("F8").Value = 523
    $ws.Range("F9").Value = 3662
    $ws.Range("F10").Value = 62
}
